$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column N (bfsGemeinde... shifts right to make room
# for the new "Traegergemeinde" column). This also grows the used range /
# <cols> width groups / dimension exactly like Excel does when a whole column
# is inserted.
$ws.Columns("N").Insert()

# Column M ("Gemeinde") becomes "Standortgemeinde" - same underlying data as
# the "Ort" column, just re-labelled / re-purposed.
$ws.Cells.Item(4, 13).Value2 = "{standortgemeindeTitle}"
$ws.Cells.Item(5, 13).Value2 = "{ort}"

# The newly inserted column N becomes "Traegergemeinde".
$ws.Cells.Item(4, 14).Value2 = "{traegergemeindeTitle}"
$ws.Cells.Item(5, 14).Value2 = "{traegergemeinde}"

# Reflect the author's final view/selection state.
$ws.Application.ActiveWindow.ScrollColumn = 7
$sel = $ws.Range("M16")
$sel.Select()
